{"js": "// Lit Review - Part 2\n//\n// Bold the \"Kaur, H., Nori, H., Jenkins, S., Caruana, R., Wallach, H., &\n// Wortman Vaughan, J. (2020).\" bibliography entry (all three runs, plus\n// the paragraph mark itself), matching the target OOXML exactly:\n// every run's <w:rPr> - and the paragraph's <w:pPr><w:rPr> (the mark) -\n// gain both <w:b/> and <w:bCs/>.\n\nconst body = context.document.body;\n\n// Locate the target bibliography paragraph via its distinctive text.\nconst results = body.search(\"Kaur, H., Nori, H., Jenkins, S.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the Kaur et al. (2020) bibliography entry.\");\n}\n\nconst para = results.items[0].paragraphs.getFirst();\n\n// Step 1: turn on bold using the native Word API. This correctly bolds\n// the paragraph mark (pPr/rPr) as well as every run in the paragraph\n// (adding <w:b/> in all four spots), but the Word JS API has no\n// complex-script-bold (\"bCs\") property to set alongside it.\npara.font.bold = true;\nawait context.sync();\n\n// Step 2: fetch the now-updated paragraph OOXML and add the missing\n// <w:bCs/> immediately after every <w:b/> we just produced, then write\n// the paragraph back via insertOoxml. This is a minimal, surgical patch\n// that reuses the engine's own (already-correct) output rather than\n// reconstructing the paragraph text by hand.\nconst paraRange = para.getRange(\"Whole\");\nconst ooxmlResult = para.getOoxml();\nawait context.sync();\n\nconst full = ooxmlResult.value;\nconst pStart = full.indexOf(\"<w:p \");\nconst pEnd = full.indexOf(\"</w:p>\", pStart) + \"</w:p>\".length;\nif (pStart === -1 || pEnd === -1) {\n  throw new Error(\"Unexpected OOXML shape returned for the paragraph.\");\n}\nlet pXml = full.substring(pStart, pEnd);\n\nconst bCount = (pXml.match(/<w:b\\/>/g) || []).length;\nif (bCount !== 4) {\n  throw new Error(\n    \"Expected 4 <w:b/> occurrences (paragraph mark + 3 runs) after bolding, found \" +\n      bCount +\n      \". Aborting to avoid corrupting content.\"\n  );\n}\npXml = pXml.split(\"<w:b/>\").join(\"<w:b/><w:bCs/>\");\n\nconst packageXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:body>\" +\n  pXml +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nparaRange.insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Lit Review - Part 2\n#\n# Bold the \"Kaur, H., Nori, H., Jenkins, S., Caruana, R., Wallach, H., &\n# Wortman Vaughan, J. (2020).\" bibliography entry (all three runs, plus\n# the paragraph mark itself), matching the target OOXML exactly:\n# every run's <w:rPr> - and the paragraph's <w:pPr><w:rPr> (the mark) -\n# gain both <w:b/> and <w:bCs/>.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Kaur, H., Nori, H., Jenkins, S.*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the Kaur et al. (2020) bibliography entry.\"\n}\n\n# Step 1: turn on bold (and complex-script bold) using the native Word\n# object model. Font.Bold correctly marks the paragraph mark (pPr/rPr)\n# as well as every run in the paragraph; Font.BoldBi correctly adds\n# <w:bCs/> to every run - but the paragraph mark's <w:bCs/> is not\n# mirrored by the engine, so that one spot needs a follow-up patch.\n$r = $target.Range\n$r.Font.Bold = 1\n$r.Font.BoldBi = 1\n\n# Step 2: fetch the now-updated paragraph OOXML and add the missing\n# <w:bCs/> to the paragraph mark's rPr only - a minimal, surgical patch\n# that reuses the engine's own (already-correct) output for everything\n# else rather than reconstructing the paragraph text by hand.\n$r2 = $target.Range\n$xml = $r2.WordOpenXML\n\n$pMatch = [regex]::Match($xml, '<w:p [^>]*w14:paraId=\"7E7A49D0\"[^>]*>.*?</w:p>', [System.Text.RegularExpressions.RegexOptions]::Singleline)\nif (-not $pMatch.Success) {\n    throw \"Could not locate the paragraph's XML to patch.\"\n}\n$pXml = $pMatch.Value\n\n$markOld = '<w:ind w:left=\"720\" w:hanging=\"720\"/><w:rPr><w:b/></w:rPr></w:pPr>'\n$markNew = '<w:ind w:left=\"720\" w:hanging=\"720\"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>'\nif (-not $pXml.Contains($markOld)) {\n    throw \"Paragraph mark rPr not in the expected shape; aborting to avoid corrupting content.\"\n}\n$pXml = $pXml.Replace($markOld, $markNew)\n\n$pkg = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' + $pXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r3 = $target.Range\n$r3.InsertXML($pkg)\n"}
